# Timesheet update: add entries for 20-02-2020 (Feb 20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Templates to copy formatting from (existing rows already styled correctly) ---
# Row 162 is a blank "day separator" row (style pattern s=5 / s=6 / s=5)
# Row 170 is a normal data row (style pattern s=1 / s=3 / s=1)
# Row 171 is a data row whose Task cell wraps/taller (style pattern s=1 / s=2 / s=1, custom row height)
$blankTemplate  = "A162:C162"
$normalTemplate = "A170:C170"
$tallTemplate   = "A171:C171"

# New data for Feb 20, 2020 (row number => Timestamp, Task, Location)
$entries = @(
    @{ Row = 172; Type = "blank" },
    @{ Row = 173; Type = "normal"; A = "Feb 20 10:00 to 11:00"; B = "Removed outliers from data, modified some of the code."; C = "Infimetrics" },
    @{ Row = 174; Type = "normal"; A = "Feb 20 11:00 to 12:00"; B = "Found lags in data"; C = "Infimetrics" },
    @{ Row = 175; Type = "normal"; A = "Feb 20 12:00 to 13:00"; B = "Removed lags from data"; C = "Infimetrics" },
    @{ Row = 176; Type = "normal"; A = "Feb 20 13:00 to 13:30"; B = "Ploted trend, looking for other methods to plot trend."; C = "Infimetrics" },
    @{ Row = 177; Type = "normal"; A = "Feb 20 13:30 to 14:30"; B = "Lunch"; C = "Infimetrics" },
    @{ Row = 178; Type = "normal"; A = "Feb 20 14:30 to 15:00"; B = "Concatinated data"; C = "Infimetrics" },
    @{ Row = 179; Type = "normal"; A = "Feb 20 15:00 to 16:00"; B = "Created new feature hour, dropped hourly timestamp. Started model building"; C = "Infimetrics" },
    @{ Row = 180; Type = "normal"; A = "Feb 20 16:00 to 17:00"; B = "Model builded, applied linear regresion"; C = "Infimetrics" },
    @{ Row = 181; Type = "normal"; A = "Feb 20 17:00 to 18:00"; B = "Checked accuracy of model"; C = "Infimetrics" },
    @{ Row = 182; Type = "tall";   A = "Feb 20 18:00 to 19:00"; B = "Modified code, made separate modules inside ML pipeline package, working on multiple`nmodel processing."; C = "Infimetrics" }
)

foreach ($entry in $entries) {
    $row = $entry.Row
    $destRange = "A$($row):C$($row)"

    switch ($entry.Type) {
        "blank" {
            $ws.Range($blankTemplate).Copy()
        }
        "tall" {
            $ws.Range($tallTemplate).Copy()
        }
        default {
            $ws.Range($normalTemplate).Copy()
        }
    }
    $ws.Range($destRange).PasteSpecial(-4122) | Out-Null

    if ($entry.Type -ne "blank") {
        $ws.Cells.Item($row, 1).Value2 = $entry.A
        $ws.Cells.Item($row, 2).Value2 = $entry.B
        $ws.Cells.Item($row, 3).Value2 = $entry.C
    }

    if ($entry.Type -eq "tall") {
        # The wrapped two-line task description needs the taller row height
        # (matches the height Excel auto-computed for the author's entry).
        $ws.Rows.Item($row).RowHeight = 45
    }
}

$excel.CutCopyMode = 0

# Move selection/view to the newly added last row, mirroring the author's saved view state
$ws.Range("D182").Select()
